# Add 2022-Q4 data:
#  1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before "2022-Q3").
#  2. Populate it with the fund-holding detail rows for 2022-Q4.
#  3. Insert a new row into the "总计" (summary) sheet for 2022-Q4, pushing the
#     existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# --- 1 & 2: new "2022-Q4" worksheet ------------------------------------------------
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2 - 009384
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'009384"
$q4.Range("B2").ClearFormats()
$q4.Range("C2").Value = "摩根士丹利华鑫MSCI中国A股指数增强A"
$q4.Range("D2").Value = "'0.38"
$q4.Range("D2").ClearFormats()
$q4.Range("E2").Value = "'89.37"
$q4.Range("E2").ClearFormats()
$q4.Range("F2").Value = "'1.08"
$q4.Range("F2").ClearFormats()
$q4.Range("G2").Value = "'0.0041"
$q4.Range("G2").ClearFormats()
$q4.Range("H2").Value = 9

# Row 3 - 014866
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'014866"
$q4.Range("B3").ClearFormats()
$q4.Range("C3").Value = "摩根士丹利华鑫MSCI中国A股指数增强C"
$q4.Range("D3").Value = "'0.00"
$q4.Range("D3").ClearFormats()
$q4.Range("E3").Value = "'89.37"
$q4.Range("E3").ClearFormats()
$q4.Range("F3").Value = "'1.08"
$q4.Range("F3").ClearFormats()
$q4.Range("G3").Value = 0
$q4.Range("H3").Value = 9

# Apply the bordered/bold "index" style (same as used on 总计!A2 and the header
# row of the other quarter sheets) to the header row and the A column.
$total.Range("A2").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A3").PasteSpecial(-4122)
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1

# --- 3: new row in "总计" summary sheet ---------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# Column A is a running 0-based index; renumber it sequentially now that a row
# was inserted (0,1,2,...,7) rather than leaving the shifted-down old values.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
